# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" on the Overview sheet and the
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" columns on
# the per-locale (zh-cn / de-de) sheets for the c60e3de8-... source file, to
# reflect a freshly generated handback report.

$wb = $excel.ActiveWorkbook

# ----- Overview sheet -----
$overview = $wb.Worksheets.Item("Overview")
# Row 2 corresponds to c60e3de8-f0b1-463f-83d6-957c38bb26a9.md
$overview.Range("G2").Value = "2016-09-02 02:56:30"

# ----- zh-cn sheet -----
$zhcn = $wb.Worksheets.Item("zh-cn")
# Row 2 corresponds to c60e3de8-f0b1-463f-83d6-957c38bb26a9.md
$zhcn.Range("H2").Value = "2016-09-02 02:56:26"
$zhcn.Range("K2").Value = "2016-09-02 02:56:43"

# ----- de-de sheet -----
$dede = $wb.Worksheets.Item("de-de")
# Row 2 corresponds to c60e3de8-f0b1-463f-83d6-957c38bb26a9.md
$dede.Range("H2").Value = "2016-09-02 02:56:30"
$dede.Range("K2").Value = "2016-09-02 02:56:50"
